# Fix exploration time & agent_step_time formula (#13)
# Updates Avg_Agent_Step_Time (G), Avg_Experiment_Time (H),
# Std_Agent_Step_Time (M) and Std_Experiment_Time (N) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @{ G = 8.314492450000001;  H = 448.76810782;    M = 0.8883186451482461; N = 74.48029021514826 }
    3  = @{ G = 9.597411920000001;  H = 837.21242793;    M = 1.583401387045373;  N = 213.485354250579 }
    4  = @{ G = 4.224720069999999;  H = 134.25883026;    M = 0.8253275188803025; N = 41.15293653647098 }
    5  = @{ G = 4.20222397;         H = 201.20463176;    M = 0.6779754676739468; N = 56.9338361412727 }
    6  = @{ G = 1.63645943;         H = 30.86698911;     M = 0.4453971724631998; N = 12.371041847014 }
    7  = @{ G = 1.73958172;         H = 47.02998725;     M = 0.4325281065604944; N = 18.1814224254417 }
    8  = @{ G = 0.83967795;         H = 11.13863903;     M = 0.3033605046516108; N = 5.355372929625585 }
    9  = @{ G = 0.8557333499999999; H = 16.49723366;     M = 0.2666376660388178; N = 7.939670751014813 }
    10 = @{ G = 0.43494278;         H = 4.3914748;       M = 0.1888859005646099; N = 2.719136353838094 }
    11 = @{ G = 0.46731378;         H = 7.16046904;      M = 0.1801239939444662; N = 4.705001554946197 }
    12 = @{ G = 0.26367343;         H = 2.0947919;       M = 0.1370732062826857; N = 1.582938666465755 }
    13 = @{ G = 0.2694337000000001; H = 3.40033725;      M = 0.1070338406259868; N = 2.324451434122841 }
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
}
